# Weekly refresh of "Fruta / hortaliza" data: the records (rows 2-16) get
# reassigned to different dates/weeks. Column D (Fecha), J (Volumen),
# K (Precio minimo), L (Precio maximo), M (Precio promedio ponderado),
# O (Origen) and P (Precio $/Kg) are redistributed across the rows while
# every other column (A,B,C,E,F,G,H,I,N,Q,R) stays constant for every row
# in this subset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that get shuffled between rows.
$colD = 4
$colJ = 10
$colK = 11
$colL = 12
$colM = 13
$colO = 15
$colP = 16

$firstRow = 2
$lastRow = 16

# Snapshot the current ("before") values for every row first, so the
# subsequent writes don't clobber data that still needs to be read.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $snapshot[$r] = @{
        D = $ws.Cells.Item($r, $colD).Value2()
        J = $ws.Cells.Item($r, $colJ).Value2()
        K = $ws.Cells.Item($r, $colK).Value2()
        L = $ws.Cells.Item($r, $colL).Value2()
        M = $ws.Cells.Item($r, $colM).Value2()
        O = $ws.Cells.Item($r, $colO).Value2()
        P = $ws.Cells.Item($r, $colP).Value2()
    }
}

# Maps each destination row to the source row whose values it should
# receive (row 3 is left untouched).
$rowMap = @{
    2  = 13
    3  = 3
    4  = 7
    5  = 12
    6  = 16
    7  = 5
    8  = 10
    9  = 2
    10 = 14
    11 = 4
    12 = 15
    13 = 9
    14 = 6
    15 = 8
    16 = 11
}

foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $data = $snapshot[$srcRow]

    $ws.Cells.Item($destRow, $colD).Value2 = $data.D
    $ws.Cells.Item($destRow, $colJ).Value2 = $data.J
    $ws.Cells.Item($destRow, $colK).Value2 = $data.K
    $ws.Cells.Item($destRow, $colL).Value2 = $data.L
    $ws.Cells.Item($destRow, $colM).Value2 = $data.M
    $ws.Cells.Item($destRow, $colO).Value2 = $data.O
    $ws.Cells.Item($destRow, $colP).Value2 = $data.P
}
